$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (s="1", grey fill) from A1 onto the new header
# cells F1:M1 before filling in their values, so the new header cells
# reuse the existing style index instead of creating a new one.
$ws.Range("A1").Copy()
$ws.Range("F1:M1").PasteSpecial(-4122)  # xlPasteFormats

# --- second "COMMAND/val1/val2/val3/val4" block (columns F:J) ---
$ws.Range("F1").Value = "COMMAND"
$ws.Range("G1").Value = "val1"
$ws.Range("H1").Value = "val2"
$ws.Range("I1").Value = "val3"
$ws.Range("J1").Value = "val4"

$ws.Range("F2").Value = "test2"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 4

$ws.Range("F3").Value = "test"
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 13
$ws.Range("J3").Value = 27

$ws.Range("F4").Value = "test3"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 40

# --- third block (columns K:M) with new shared strings A/B/C/D ---
$ws.Range("K1").Value = "COMMAND"
$ws.Range("L1").Value = "val1"
$ws.Range("M1").Value = "val2"

$ws.Range("K2").Value = "test"
$ws.Range("L2").Value = "A"
$ws.Range("M2").Value = "B"

$ws.Range("K3").Value = "test"
$ws.Range("L3").Value = "C"
$ws.Range("M3").Value = "D"

# Final selection ends up on K4, matching the saved workbook's view state.
$ws.Range("K4").Select()
